$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.811.94"
$ws.Range("E2").Value = "  +0.18%  "
$ws.Range("D3").Value = "2.293.32"
$ws.Range("E3").Value = "  -0.12%  "
$ws.Range("E4").Value = "  +0.16%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "116.70"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +14.45%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "268.95"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.39%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.628"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +1.47%  "
$ws.Range("E8").Value = "  +0.22%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.618"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +1.64%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "49.51"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +9.14%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0947"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +1.33%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "8.93"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +11.46%  "
$ws.Range("E13").Value = "  +0.61%  "
$ws.Range("E14").Value = "  +0.77%  "
$ws.Range("D15").Value = "2.638.06"
$ws.Range("E15").Value = "  -0.13%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.882"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +2.78%  "
$ws.Range("D17").Value = "2.298.09"
$ws.Range("E17").Value = "  -0.75%  "
$ws.Range("D18").Value = "43.704.59"
$ws.Range("E18").Value = "  -0.09%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.0000109"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -0.74%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "6.99"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +12.00%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "72.45"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +0.27%  "
$ws.Range("E22").Value = "  -1.45%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "9.96"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +8.09%  "
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("E25").Value = "  +3.57%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "11.72"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +4.81%  "
$ws.Range("E28").Value = "  +0.78%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "42.00"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +6.32%  "
$ws.Range("E30").Value = "  -2.01%  "
$ws.Range("E31").Value = "  -2.03%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "173.44"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -2.31%  "
$ws.Range("E33").Value = "  +4.47%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "21.64"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -0.71%  "
$ws.Range("E35").Value = "  +5.09%  "
$ws.Range("E36").Value = "  +0.60%  "
$ws.Range("E37").Value = "  -1.90%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.0359"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +1.74%  "
$ws.Range("E39").Value = "  +0.28%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "3.81"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +7.16%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "14.64"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +18.88%  "
$ws.Range("E42").Value = "  +4.32%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "74.12"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +14.35%  "
$ws.Range("E44").Value = "  +3.44%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "6.40"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +22.19%  "
$ws.Range("E46").Value = "  +0.19%  "
$ws.Range("E47").Value = "  +0.53%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "8.75"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -1.08%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "103.29"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +5.21%  "
$ws.Range("E50").Value = "  +4.36%  "
$ws.Range("E51").Value = "  -1.41%  "
